# A1.4.xlsx — "Add files via upload" re-save.
#
# The canonical diff is dominated by cosmetic re-save noise (Excel build/
# rupBuild numbers, xr/xr2/xr3 revision-tracking namespaces, the author's
# per-machine x15ac:absPath, window pixel geometry, x14ac:dyDescent hints,
# a different default-row-height baseline, etc.) that simply falls out of
# opening the file in a different Excel build/profile and saving again —
# it is not something a user deliberately edits through the object model,
# so there is nothing meaningful to script for it here.
#
# The deliberate, content-level edit visible in the diff is on the
# "#A1.4" worksheet: the worked "absolute vs relative reference" example
# had its helper formulas in column C (rows 3-5) and column G (rows 2-5)
# cleared, leaving only the original C2 formula (=A2+B2) behind.

$wb = $excel.ActiveWorkbook

$wsExample = $wb.Worksheets.Item("#A1.4")

# Clear the shared "C" formulas in rows 3-5 (C2's "=A2+B2" is left alone).
$wsExample.Range("C3:C5").ClearContents()

# Clear every "G" formula (rows 2-5).
$wsExample.Range("G2:G5").ClearContents()

# Turn on iterative calculation (calcPr iterate="1" in the saved file).
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.001

# The workbook's active tab ends up on "A1.4" (activeTab moves from the
# "#A1.4" example sheet to the adjoining "A1.4" sheet).
$wsFinal = $wb.Worksheets.Item("A1.4")
$wsFinal.Activate()
